$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the subtitle placeholder shape ("Subtitle 2", a ppPlaceholderSubtitle
# placeholder) from slide 1. Because the slide layout still defines this
# placeholder, deleting it the first time only resets it back to its empty
# layout default (it stays in Shapes, just emptied out and renumbered) -
# deleting it a second time fully removes it from the slide.
for ($pass = 1; $pass -le 2; $pass++) {
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shape = $s.Shapes.Item($i)
        if ($shape.Type -eq 14 -and $shape.PlaceholderFormat.Type -eq 4) {
            $shape.Delete()
        }
    }
}
